$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 'ba'
$ws.Range("J2").Value = 'Appreciation'
$ws.Range("I3").Value = 'aa'
$ws.Range("J3").Value = 'Agree/Accept'
$ws.Range("I5").Value = 'b'
$ws.Range("J5").Value = 'Acknowledge (Backchannel)'
$ws.Range("I7").Value = 'sd'
$ws.Range("J7").Value = 'Statement-non-opinion'
$ws.Range("I12").Value = 'aa'
$ws.Range("J12").Value = 'Agree/Accept'
$ws.Range("I16").Value = 'b'
$ws.Range("J16").Value = 'Acknowledge (Backchannel)'
$ws.Range("I17").Value = 'b'
$ws.Range("J17").Value = 'Acknowledge (Backchannel)'
$ws.Range("I18").Value = 'b'
$ws.Range("J18").Value = 'Acknowledge (Backchannel)'
$ws.Range("I19").Value = 'aa'
$ws.Range("J19").Value = 'Agree/Accept'
$ws.Range("I25").Value = '%'
$ws.Range("J25").Value = 'Uninterpretable'
$ws.Range("I31").Value = 'aa'
$ws.Range("J31").Value = 'Agree/Accept'
$ws.Range("I37").Value = '%'
$ws.Range("J37").Value = 'Uninterpretable'
$ws.Range("I46").Value = 'aa'
$ws.Range("J46").Value = 'Agree/Accept'
$ws.Range("I49").Value = 'sd'
$ws.Range("J49").Value = 'Statement-non-opinion'
$ws.Range("I73").Value = 'sv'
$ws.Range("J73").Value = 'Statement-opinion'
$ws.Range("I75").Value = 'aa'
$ws.Range("J75").Value = 'Agree/Accept'
$ws.Range("I78").Value = 'aa'
$ws.Range("J78").Value = 'Agree/Accept'
$ws.Range("I81").Value = 'ba'
$ws.Range("J81").Value = 'Appreciation'
$ws.Range("I82").Value = 'sv'
$ws.Range("J82").Value = 'Statement-opinion'
$ws.Range("I97").Value = 'aa'
$ws.Range("J97").Value = 'Agree/Accept'
$ws.Range("I105").Value = '%'
$ws.Range("J105").Value = 'Uninterpretable'
$ws.Range("I110").Value = '%'
$ws.Range("J110").Value = 'Uninterpretable'
$ws.Range("I111").Value = '%'
$ws.Range("J111").Value = 'Uninterpretable'
$ws.Range("I112").Value = 'sd'
$ws.Range("J112").Value = 'Statement-non-opinion'
$ws.Range("I116").Value = 'sd'
$ws.Range("J116").Value = 'Statement-non-opinion'
$ws.Range("I131").Value = 'sv'
$ws.Range("J131").Value = 'Statement-opinion'
$ws.Range("I142").Value = '%'
$ws.Range("J142").Value = 'Uninterpretable'
$ws.Range("I143").Value = '%'
$ws.Range("J143").Value = 'Uninterpretable'
$ws.Range("I153").Value = 'sv'
$ws.Range("J153").Value = 'Statement-opinion'
$ws.Range("I155").Value = 'aa'
$ws.Range("J155").Value = 'Agree/Accept'
$ws.Range("I156").Value = 'aa'
$ws.Range("J156").Value = 'Agree/Accept'
$ws.Range("I157").Value = 'aa'
$ws.Range("J157").Value = 'Agree/Accept'
$ws.Range("I158").Value = 'aa'
$ws.Range("J158").Value = 'Agree/Accept'
$ws.Range("I160").Value = '%'
$ws.Range("J160").Value = 'Uninterpretable'
$ws.Range("I179").Value = 'aa'
$ws.Range("J179").Value = 'Agree/Accept'
$ws.Range("I183").Value = 'sd'
$ws.Range("J183").Value = 'Statement-non-opinion'
$ws.Range("I188").Value = '%'
$ws.Range("J188").Value = 'Uninterpretable'
$ws.Range("I189").Value = '%'
$ws.Range("J189").Value = 'Uninterpretable'
